$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting existing rows 12-61 down to 13-62.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with its data.
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(12, 3).Value = "Bíobío"
$ws.Cells.Item(12, 4).Value = 45238
$ws.Cells.Item(12, 5).Value = 8
$ws.Cells.Item(12, 6).Value = "Fruta"
$ws.Cells.Item(12, 7).Value = 100107
$ws.Cells.Item(12, 8).Value = "Otros"
$ws.Cells.Item(12, 9).Value = 100107002
$ws.Cells.Item(12, 10).Value = "Chirimoya"
$ws.Cells.Item(12, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(12, 12).Value = "Primera"
$ws.Cells.Item(12, 13).Value = 100
$ws.Cells.Item(12, 14).Value = 19000
$ws.Cells.Item(12, 15).Value = 19000
$ws.Cells.Item(12, 16).Value = 19000
$ws.Cells.Item(12, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(12, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(12, 19).Value = 1900
$ws.Cells.Item(12, 20).Value = 10
